$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextRow($row, $values) {
    # $values is an ordered list of text to place starting at column C (3)
    # through F (6). Excel would otherwise auto-detect numeric-looking
    # strings like "24.000" or "0.56" and coerce them to numbers, so we
    # temporarily force a text number format, write the values, then reset
    # the number format back to Normal so the cells keep the workbook's
    # default (unstyled) appearance, matching the source data which stores
    # these as plain inline strings with no special style.
    $rng = $ws.Range("C" + $row + ":F" + $row)
    $rng.NumberFormat = "@"
    $col = 3
    foreach ($v in $values) {
        $ws.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
    $rng.Style = "Normal"
}

# --- Update the date in B2 ---
$ws.Cells.Item(2, 2).Value = "22-06-24"

# --- Row 3: Party ---
$ws.Cells.Item(3, 1).Value = "Party : "
$ws.Cells.Item(3, 2).Value = "hello"
$ws.Range("B3:F3").Merge()

# --- Row 4: Job Name ---
$ws.Cells.Item(4, 1).Value = "Job Name : "
$ws.Cells.Item(4, 2).Value = "WRP MGX KRE ROU 8TW 29G+11.8G(4008652)"
$ws.Range("B4:F4").Merge()

# --- Row 5: table header (1st table), shaded light-gray ---
$ws.Cells.Item(5, 1).Value = "SNo."
$ws.Cells.Item(5, 2).Value = "Roll No."
$ws.Cells.Item(5, 3).Value = "Gross Wt."
$ws.Cells.Item(5, 4).Value = "Tare Wt."
$ws.Cells.Item(5, 5).Value = "Core Wt."
$ws.Cells.Item(5, 6).Value = "Net Wt."
$ws.Range("A5:F5").Interior.Color = 13882323

# --- Row 6 (roll 15, status Done) ---
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = 15
Set-TextRow 6 @("24.000", "0.56", "0.2", "23.240")
$ws.Cells.Item(6, 7).Value = "Done"

# --- Row 7 (roll 16, status Done) ---
$ws.Cells.Item(7, 1).Value = 2
$ws.Cells.Item(7, 2).Value = 16
Set-TextRow 7 @("15.000", "0.56", "0.2", "14.240")
$ws.Cells.Item(7, 7).Value = "Done"

# --- Row 8 (roll 17, status Pending) ---
$ws.Cells.Item(8, 1).Value = 3
$ws.Cells.Item(8, 2).Value = 17
Set-TextRow 8 @("20.000", "0.56", "0.2", "19.240")
$ws.Cells.Item(8, 7).Value = "Pending"

# --- Row 9: totals for table 1 ---
$ws.Cells.Item(9, 2).Value = "Total"
$ws.Cells.Item(9, 3).Value = 59
$ws.Cells.Item(9, 4).Value = 1.68
$ws.Cells.Item(9, 5).Value = 0.6
$ws.Cells.Item(9, 6).Value = 56.72

# --- Row 10: second Job Name ---
$ws.Cells.Item(10, 1).Value = "Job Name : "
$ws.Cells.Item(10, 2).Value = "bad vaiety"
$ws.Range("B10:F10").Merge()

# --- Row 11: table header (2nd table), shaded light-gray, extends to G ---
$ws.Cells.Item(11, 1).Value = "SNo."
$ws.Cells.Item(11, 2).Value = "Roll No."
$ws.Cells.Item(11, 3).Value = "Gross Wt."
$ws.Cells.Item(11, 4).Value = "Tare Wt."
$ws.Cells.Item(11, 5).Value = "Core Wt."
$ws.Cells.Item(11, 6).Value = "Net Wt."
$ws.Range("A11:G11").Interior.Color = 13882323

# --- Row 12 (roll 18, status Pending) ---
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = 18
Set-TextRow 12 @("20.000", "0.56", "0.2", "19.240")
$ws.Cells.Item(12, 7).Value = "Pending"

# --- Row 13: totals for table 2 ---
$ws.Cells.Item(13, 2).Value = "Total"
$ws.Cells.Item(13, 3).Value = 20
$ws.Cells.Item(13, 4).Value = 0.5600000000000001
$ws.Cells.Item(13, 5).Value = 0.2
$ws.Cells.Item(13, 6).Value = 19.24
